# Chưa xử lý chấm công cho toàn bộ data
# Adds a new "Sheet2" (placed after the existing "Sheet1") containing a
# 22-row x 19-column numeric grid, and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so the sheet order is Sheet1, Sheet2
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Build the A1:S22 numeric grid.
# Column A = row number.
# Even columns (B, D, F, H, J, L, N, P, R) only have a value on odd rows,
#   and that value is constant, equal to the column number (B=2, D=4, ...).
# Column C and all further odd columns (C, E, G, I, K, M, O, Q, S) have a
#   value on every row that decreases linearly as the row grows:
#   value(col, row) = col - (col - 2) * (row - 1)
$data = New-Object 'object[,]' 22,19
for ($r = 1; $r -le 22; $r++) {
    for ($ci = 1; $ci -le 19; $ci++) {
        if ($ci -eq 1) {
            $val = $r
        } elseif ($ci % 2 -eq 0) {
            if ($r % 2 -eq 1) { $val = $ci } else { $val = $null }
        } else {
            $slope = -($ci - 2)
            $val = $ci + $slope * ($r - 1)
        }
        $data[$r - 1, $ci - 1] = $val
    }
}
$ws2.Range("A1:S22").Value = $data

# Match the saved selection/active cell on the new sheet
$ws2.Range("P11").Select() | Out-Null

# Make Sheet2 the active/visible tab
$ws2.Activate() | Out-Null
